$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "Datos actualizados ..." timestamp refresh
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 01:46"

# Rows with pure numeric refresh (country/ranking position unchanged)
$ws.Range("B9").Value = 19469
$ws.Range("C9").Value = 5680
$ws.Range("D9").Value = 147
$ws.Range("E9").Value = 19064
$ws.Range("F9").Value = 64
$ws.Range("G9").Value = 51
$ws.Range("H9").Value = 258

$ws.Range("B12").Value = 5615
$ws.Range("C12").Value = 1393
$ws.Range("D12").Value = 15
$ws.Range("E12").Value = 5544
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 13
$ws.Range("H12").Value = 56

$ws.Range("B20").Value = 1087
$ws.Range("C20").Value = 214
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = 1061
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 12

$ws.Range("B25").Value = 928
$ws.Range("C25").Value = 172
$ws.Range("D25").Value = 46
$ws.Range("E25").Value = 875
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 7

# Rows 76-79: Marruecos case counts rose, moving it up past Jordania/Hungria/Islas Feroe
$ws.Range("A76").Value = "Marruecos"
$ws.Range("B76").Value = 86
$ws.Range("C76").Value = 23
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 81
$ws.Range("F76").Value = 1
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 3

$ws.Range("A77").Value = "Jordania"
$ws.Range("B77").Value = 85
$ws.Range("C77").Value = 16
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 84
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0

$ws.Range("A78").Value = "Hungria"
$ws.Range("B78").Value = 85
$ws.Range("C78").Value = 12
$ws.Range("D78").Value = 7
$ws.Range("E78").Value = 74
$ws.Range("F78").Value = 6
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 4

$ws.Range("A79").Value = "Islas Feroe"
$ws.Range("B79").Value = 80
$ws.Range("C79").Value = 8
$ws.Range("D79").Value = 3
$ws.Range("E79").Value = 77
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 0

# Rows 110-119: Cuba case counts rose, moving it up past Camerun..Bolivia
$ws.Range("A110").Value = "Cuba"
$ws.Range("B110").Value = 21
$ws.Range("C110").Value = 10
$ws.Range("D110").Value = 0
$ws.Range("E110").Value = 20
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 1

$ws.Range("A111").Value = "Camerun"
$ws.Range("B111").Value = 20
$ws.Range("C111").Value = 7
$ws.Range("D111").Value = 2
$ws.Range("E111").Value = 18
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 0

$ws.Range("A112").Value = "Banglades"
$ws.Range("B112").Value = 20
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 3
$ws.Range("E112").Value = 16
$ws.Range("F112").Value = 1
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 1

$ws.Range("A113").Value = "Jamaica"
$ws.Range("B113").Value = 19
$ws.Range("C113").Value = 4
$ws.Range("D113").Value = 2
$ws.Range("E113").Value = 16
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1

$ws.Range("A114").Value = "Paraguay"
$ws.Range("B114").Value = 18
$ws.Range("C114").Value = 5
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 18
$ws.Range("F114").Value = 1
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 0

$ws.Range("A115").Value = "Consejo Danes para los Refugiados"
$ws.Range("B115").Value = 18
$ws.Range("C115").Value = 4
$ws.Range("D115").Value = 0
$ws.Range("E115").Value = 18
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

$ws.Range("A116").Value = "Ruanda"
$ws.Range("B116").Value = 17
$ws.Range("C116").Value = 6
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 17
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

$ws.Range("A117").Value = "Macao"
$ws.Range("B117").Value = 17
$ws.Range("C117").Value = 0
$ws.Range("D117").Value = 10
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0

$ws.Range("A118").Value = "Ghana"
$ws.Range("B118").Value = 16
$ws.Range("C118").Value = 5
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 16
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 0

$ws.Range("A119").Value = "Bolivia"
$ws.Range("B119").Value = 16
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 16
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0
